$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.90"
$ws.Range("E2").Value = "1BNBBNBBestin24h"

$ws.Range("D4").Value = "'5.555"

$ws.Range("D5").Value = "'0.05637"

$ws.Range("D6").Value = "'3.403"

$ws.Range("D7").Value = "'6.480"

$ws.Range("D8").Value = "'0.8010"

$ws.Range("D9").Value = "'1.070"

$ws.Range("D10").Value = "'0.1430"

$ws.Range("D11").Value = "'0.07408"

$ws.Range("D12").Value = "'0.03190"

$ws.Range("D13").Value = "'0.02970"

$ws.Range("D14").Value = "'0.09261"

$ws.Range("D15").Value = "'0.001667"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'2.981"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04711"
$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005741"
$ws.Range("E18").Value = "17OneONE"

$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "'0.006261"
$ws.Range("E19").Value = "18TigerCashTCH"

$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Value = "'0.001057"
$ws.Range("E20").Value = "19BitKanKAN"

$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").Value = "'0.003823"
$ws.Range("E21").Value = "20HotbitTokenHTB"

$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "'0.0001501"
$ws.Range("E22").Value = "21NitroExNTX"

$ws.Range("B23").Value = "UpBots"
$ws.Range("C23").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D23").Value = "'0.0004601"
$ws.Range("E23").Value = "22UpBotsUBXT"

$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").Value = "'3.983"
$ws.Range("E24").Value = "23LEOLEO"

$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").Value = "'2.112"
$ws.Range("E25").Value = "24BTSETokenBTSE"

$ws.Range("B26").Value = "BitpandaEcosystemToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D26").Value = "'0.3311"
$ws.Range("E26").Value = "25BitpandaEcosystemTokenBEST"

$ws.Range("B27").Value = "ProBitToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D27").Value = "'0.1277"
$ws.Range("E27").Value = "26ProBitTokenPROB"

$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1046"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("D42").Value = "'0.002971"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003248"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

$ws.Range("D44").Value = "'0.008678"

$ws.Range("D45").Value = "'0.00005631"

$ws.Range("D47").Value = "'0.6802"

$ws.Range("D48").Value = "'0.02797"

$ws.Range("D49").Value = "'0.00002100"
